$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'71.677.10"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "'3.814.78"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'705.31"
$ws.Range("E5").Value = "  +6.40%  "
$ws.Range("D6").Value = "'175.11"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("D7").Value = "'3.812.88"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'7.41"
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("D14").Value = "'36.66"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").Value = "'4.460.09"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "'3.797.48"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'71.617.32"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "'17.78"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'7.25"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'0.115"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("E21").Value = "  +6.61%  "
$ws.Range("D22").Value = "'484.98"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "'84.69"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "'12.35"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'10.59"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "'3.966.63"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +10.89%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'7.62"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").Value = "'0.187"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("D35").Value = "'29.67"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "'9.33"
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").Value = "'3.48"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "'2.36"
$ws.Range("E40").Value = "  +12.23%  "
$ws.Range("D41").Value = "'6.06"
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").Value = "'0.995"
$ws.Range("E42").Value = "  +2.73%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'164.82"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").Value = "'0.000309"
$ws.Range("E46").Value = "  +8.42%  "
$ws.Range("D47").Value = "'44.84"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'48.70"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "'421.25"
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'1.39"
$ws.Range("E51").Value = "  -3.24%  "
